# Fruta / hortaliza, semanal
# Insert two new weekly price rows (2023-08-09) for "Piña" / "Femacal de La Calera"
# right before the existing row 1063, pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 1063:1064 - this shifts rows 1063-1107 down to 1065-1109
# and keeps the date-formatted style (s=2) of column D which Excel copies from the
# row above during an insert.
$ws.Rows("1063:1064").Insert()

# New row 1063: "Primera" quality, 108 units, $22.000, $/caja 12 unidades
$ws.Range("A1063").Value = 3
$ws.Range("B1063").Value = "Femacal de La Calera"
$ws.Range("C1063").Value = "Coquimbo"
$ws.Range("D1063").Value = 45147
$ws.Range("E1063").Value = 5
$ws.Range("F1063").Value = "Fruta"
$ws.Range("G1063").Value = 100108
$ws.Range("H1063").Value = "Tropicales y subtropicales"
$ws.Range("I1063").Value = 100108005
$ws.Range("J1063").Value = "Piña"
$ws.Range("K1063").Value = "Caramelo"
$ws.Range("L1063").Value = "Primera"
$ws.Range("M1063").Value = 108
$ws.Range("N1063").Value = 22000
$ws.Range("O1063").Value = 22000
$ws.Range("P1063").Value = 22000
$ws.Range("Q1063").Value = "`$/caja 12 unidades"
$ws.Range("R1063").Value = "Ecuador"
$ws.Range("S1063").Value = 1833
$ws.Range("T1063").Value = 12

# New row 1064: "Segunda" quality, 108 units, $22.000, $/caja 14 unidades
$ws.Range("A1064").Value = 3
$ws.Range("B1064").Value = "Femacal de La Calera"
$ws.Range("C1064").Value = "Coquimbo"
$ws.Range("D1064").Value = 45147
$ws.Range("E1064").Value = 5
$ws.Range("F1064").Value = "Fruta"
$ws.Range("G1064").Value = 100108
$ws.Range("H1064").Value = "Tropicales y subtropicales"
$ws.Range("I1064").Value = 100108005
$ws.Range("J1064").Value = "Piña"
$ws.Range("K1064").Value = "Caramelo"
$ws.Range("L1064").Value = "Segunda"
$ws.Range("M1064").Value = 108
$ws.Range("N1064").Value = 22000
$ws.Range("O1064").Value = 22000
$ws.Range("P1064").Value = 22000
$ws.Range("Q1064").Value = "`$/caja 14 unidades"
$ws.Range("R1064").Value = "Ecuador"
$ws.Range("S1064").Value = 1571
$ws.Range("T1064").Value = 14

Write-Host "Inserted rows 1063-1064 and populated new weekly data."
